$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3923.423
$ws.Range("I132").Value = 4064.36
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 12193.08
$ws.Range("L132").Value = 1200
$ws.Range("M132").Value = -9663.08
$ws.Range("N132").Value = -6260
$ws.Range("H137").Value = 25003096
$ws.Range("I137").Value = 1586.0476
$ws.Range("J137").Value = 52636344
$ws.Range("K137").Value = 4758.142800000001
$ws.Range("L137").Value = 157909032
$ws.Range("M137").Value = -2208.142800000001
$ws.Range("N137").Value = -157914132
$ws.Range("H138").Value = 2728.1968
$ws.Range("I138").Value = 2270
$ws.Range("J138").Value = 3201.6667
$ws.Range("K138").Value = 6810
$ws.Range("L138").Value = 9605.000100000001
$ws.Range("M138").Value = -1670
$ws.Range("N138").Value = -19885.0001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4904499
$ws.Range("I2").Value = 3046
$ws.Range("J2").Value = 29411764
$ws.Range("K2").Value = 3046
$ws.Range("L2").Value = 29411764
$ws.Range("M2").Value = -2933
$ws.Range("N2").Value = -29411990
$ws.Range("H37").Value = 9518.5
$ws.Range("I37").Value = 2034
$ws.Range("J37").Value = 10587.714
$ws.Range("K37").Value = 2034
$ws.Range("L37").Value = 10587.714
$ws.Range("M37").Value = -1761
$ws.Range("N37").Value = -11133.714
$ws.Range("H116").Value = 4904499
$ws.Range("I116").Value = 3046
$ws.Range("J116").Value = 29411764
$ws.Range("K116").Value = 3046
$ws.Range("L116").Value = 29411764
$ws.Range("M116").Value = -752
$ws.Range("N116").Value = -29416352
$ws.Range("H132").Value = 2877639
$ws.Range("I132").Value = 3835186.2
$ws.Range("J132").Value = 4997.3335
$ws.Range("K132").Value = 11505558.6
$ws.Range("L132").Value = 14992.0005
$ws.Range("M132").Value = -11503028.6
$ws.Range("N132").Value = -20052.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4904499
$ws.Range("I3").Value = 3046
$ws.Range("J3").Value = 29411764
$ws.Range("K3").Value = 3046
$ws.Range("L3").Value = 29411764
$ws.Range("M3").Value = -2932
$ws.Range("N3").Value = -29411992
$ws.Range("H20").Value = 2784.4285
$ws.Range("I20").Value = 3526.8572
$ws.Range("J20").Value = 2042
$ws.Range("K20").Value = 3526.8572
$ws.Range("L20").Value = 2042
$ws.Range("M20").Value = -3279.8572
$ws.Range("N20").Value = -2536
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6085.5586
$ws.Range("I31").Value = 2111.4
$ws.Range("J31").Value = 7741.4585
$ws.Range("K31").Value = 2111.4
$ws.Range("L31").Value = 7741.4585
$ws.Range("M31").Value = -1816.4
$ws.Range("N31").Value = -8331.458500000001
$ws.Range("H34").Value = 6085.5586
$ws.Range("I34").Value = 2111.4
$ws.Range("J34").Value = 7741.4585
$ws.Range("K34").Value = 2111.4
$ws.Range("L34").Value = 7741.4585
$ws.Range("M34").Value = -1909.4
$ws.Range("N34").Value = -8145.4585
$ws.Range("H50").Value = 10875.2
$ws.Range("J50").Value = 10875.2
$ws.Range("L50").Value = 10875.2
$ws.Range("N50").Value = -12125.2
$ws.Range("H51").Value = 9887.111000000001
$ws.Range("J51").Value = 10197.714
$ws.Range("L51").Value = 10197.714
$ws.Range("N51").Value = -11669.714
$ws.Range("H60").Value = 9823.444
$ws.Range("J60").Value = 10051.375
$ws.Range("L60").Value = 10051.375
$ws.Range("N60").Value = -11073.375
$ws.Range("H61").Value = 9887.111000000001
$ws.Range("J61").Value = 10197.714
$ws.Range("L61").Value = 10197.714
$ws.Range("N61").Value = -10893.714
$ws.Range("H68").Value = 18406.572
$ws.Range("J68").Value = 19096.334
$ws.Range("L68").Value = 19096.334
$ws.Range("N68").Value = -20594.334
$ws.Range("H71").Value = 18406.572
$ws.Range("J71").Value = 19096.334
$ws.Range("L71").Value = 57289.00199999999
$ws.Range("N71").Value = -64777.00199999999
$ws.Range("H74").Value = 14421.2
$ws.Range("J74").Value = 16730.25
$ws.Range("L74").Value = 16730.25
$ws.Range("N74").Value = -18478.25
$ws.Range("H77").Value = 14421.2
$ws.Range("J77").Value = 16730.25
$ws.Range("L77").Value = 50190.75
$ws.Range("N77").Value = -58926.75
$ws.Range("H93").Value = 12431.4
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H140").Value = 71314.28999999999
$ws.Range("J140").Value = 79840
$ws.Range("L140").Value = 79840
$ws.Range("N140").Value = -90200
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 2694.4
$ws.Range("J100").Value = 2767.375
$ws.Range("L100").Value = 8302.125
$ws.Range("N100").Value = -9924.125
$ws.Range("H105").Value = 446001500
$ws.Range("J105").Value = 446001500
$ws.Range("L105").Value = 1338004500
$ws.Range("N105").Value = -1338009742
$ws.Range("H131").Value = 53035960
$ws.Range("J131").Value = 27781354
$ws.Range("L131").Value = 83344062
$ws.Range("N131").Value = -83354142
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 9550
$ws.Range("I59").Value = 9100
$ws.Range("J59").Value = 10000
$ws.Range("K59").Value = 9100
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = -8517
$ws.Range("N59").Value = -11166
$ws.Range("H94").Value = 30224
$ws.Range("J94").Value = 30224
$ws.Range("L94").Value = 30224
$ws.Range("N94").Value = -31576
$ws.Range("H97").Value = 1773.3334
$ws.Range("I97").Value = 1567.5
$ws.Range("J97").Value = 2596.6667
$ws.Range("K97").Value = 1567.5
$ws.Range("L97").Value = 2596.6667
$ws.Range("M97").Value = -1071.5
$ws.Range("N97").Value = -3588.6667
$ws.Range("H107").Value = 393.78262
$ws.Range("I107").Value = 337.8
$ws.Range("J107").Value = 498.75
$ws.Range("K107").Value = 337.8
$ws.Range("L107").Value = 498.75
$ws.Range("M107").Value = 1582.2
$ws.Range("N107").Value = -4338.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2256.25
$ws.Range("I68").Value = 2026.9231
$ws.Range("J68").Value = 2527.2727
$ws.Range("K68").Value = 2026.9231
$ws.Range("L68").Value = 2527.2727
$ws.Range("M68").Value = -1277.9231
$ws.Range("N68").Value = -4025.2727
$ws.Range("H71").Value = 2256.25
$ws.Range("I71").Value = 2026.9231
$ws.Range("J71").Value = 2527.2727
$ws.Range("K71").Value = 10134.6155
$ws.Range("L71").Value = 12636.3635
$ws.Range("M71").Value = -6390.6155
$ws.Range("N71").Value = -20124.3635
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H134").Value = 27429
$ws.Range("J134").Value = 27429
$ws.Range("L134").Value = 27429
$ws.Range("N134").Value = -37569
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1448573
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 3366670.2
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 3366670.2
$ws.Range("M14").Value = -9832
$ws.Range("N14").Value = -3367006.2
$ws.Range("H107").Value = 941.2308
$ws.Range("I107").Value = 624.1
$ws.Range("J107").Value = 1998.3334
$ws.Range("K107").Value = 1872.3
$ws.Range("L107").Value = 5995.0002
$ws.Range("M107").Value = 47.69999999999982
$ws.Range("N107").Value = -9835.0002
